# Update cryptocurrency price/volume data per latest symbol-list refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'275.89"
$ws.Range("E2").Value = "'0.76%"
$ws.Range("D3").Value = "'27.25"
$ws.Range("E3").Value = "'1.94%"
$ws.Range("D4").Value = "'4.860"
$ws.Range("E4").Value = "'-0.87%"
$ws.Range("D5").Value = "'0.06410"
$ws.Range("E5").Value = "'1.24%"
$ws.Range("D6").Value = "'6.925"
$ws.Range("E6").Value = "'0.36%"
$ws.Range("D7").Value = "'1.219"
$ws.Range("E7").Value = "'-1.69%"
$ws.Range("D8").Value = "'0.8798"
$ws.Range("E8").Value = "'0.13%"
$ws.Range("D9").Value = "'0.1514"
$ws.Range("E9").Value = "'3.95%"
$ws.Range("D10").Value = "'0.05085"
$ws.Range("E10").Value = "'2.43%"
$ws.Range("D11").Value = "'0.07592"
$ws.Range("E11").Value = "'3.36%"
$ws.Range("D12").Value = "'0.02981"
$ws.Range("E12").Value = "'-4.48%"
$ws.Range("D13").Value = "'0.08996"
$ws.Range("E13").Value = "'-0.68%"
$ws.Range("D14").Value = "'0.001566"
$ws.Range("E14").Value = "'-1.82%"
$ws.Range("D15").Value = "'0.0006407"
$ws.Range("E15").Value = "'1.60%"
$ws.Range("D16").Value = "'0.006188"
$ws.Range("E16").Value = "'2.71%"
$ws.Range("D17").Value = "'3.472"
$ws.Range("E17").Value = "'0.16%"
$ws.Range("D18").Value = "'3.310"
$ws.Range("E18").Value = "'-1.56%"
$ws.Range("D19").Value = "'2.284"
$ws.Range("E19").Value = "'0.57%"
$ws.Range("E20").Value = "'-0.96%"
$ws.Range("D21").Value = "'0.1362"
$ws.Range("E21").Value = "'2.53%"
$ws.Range("D22").Value = "'3.912"
$ws.Range("E22").Value = "'-0.02%"
$ws.Range("D23").Value = "'0.04435"
$ws.Range("E23").Value = "'0.37%"
$ws.Range("D24").Value = "'0.001176"
$ws.Range("E24").Value = "'-0.25%"
$ws.Range("D25").Value = "'0.004264"
$ws.Range("E25").Value = "'15.34%"
$ws.Range("D26").Value = "'0.0001202"
$ws.Range("E26").Value = "'-0.03%"
$ws.Range("E27").Value = "'13.70%"
$ws.Range("D40").Value = "'0.04155"
$ws.Range("E40").Value = "'2.64%"
$ws.Range("D41").Value = "'0.006811"
$ws.Range("E41").Value = "'2.56%"
$ws.Range("D42").Value = "'0.1174"
$ws.Range("E42").Value = "'0.64%"
$ws.Range("D43").Value = "'0.002173"
$ws.Range("E43").Value = "'3.30%"
$ws.Range("D44").Value = "'0.01187"
$ws.Range("E44").Value = "'-4.20%"
$ws.Range("D45").Value = "'0.00005174"
$ws.Range("E45").Value = "'-3.07%"
$ws.Range("D46").Value = "'1.672"
$ws.Range("E46").Value = "'-29.03%"
$ws.Range("E47").Value = "'-0.19%"
